$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.092285666666667
$ws.Range("H2").Value = 27.276857
$ws.Range("I2").Value = 0.4104940601903898
$ws.Range("J2").Value = 0.4104940601903898
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.273683333333334
$ws.Range("N2").Value = 18.82105
$ws.Range("O2").Value = 0.5180127793393899
$ws.Range("P2").Value = 0.51801277933939
$ws.Range("Q2").Value = 57.04212104887223
$ws.Range("R2").Value = 513.37908943985
$ws.Range("S2").Value = 0.2126411690215347
$ws.Range("T2").Value = 0.2126411690215347
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.092285666666667
$ws.Range("H3").Value = 27.276857
$ws.Range("I3").Value = 0.4104940601903898
$ws.Range("J3").Value = 0.4104940601903898
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.455641666666668
$ws.Range("N3").Value = 16.366925
$ws.Range("O3").Value = 0.4504677639392779
$ws.Range("P3").Value = 0.450467763939278
$ws.Range("Q3").Value = 49.60425252830279
$ws.Range("R3").Value = 446.4382727547251
$ws.Range("S3").Value = 0.1849143414043203
$ws.Range("T3").Value = 0.1849143414043203
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.092285666666667
$ws.Range("H4").Value = 27.276857
$ws.Range("I4").Value = 0.4104940601903898
$ws.Range("J4").Value = 0.4104940601903898
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.381734
$ws.Range("N4").Value = 1.145202
$ws.Range("O4").Value = 0.03151945672133213
$ws.Range("P4").Value = 0.03151945672133213
$ws.Range("Q4").Value = 3.470834576679333
$ws.Range("R4").Value = 31.237511190114
$ws.Range("S4").Value = 0.0129385497645349
$ws.Range("T4").Value = 0.01293854976453489
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.532627333333332
$ws.Range("H5").Value = 28.597882
$ws.Range("I5").Value = 0.4303743900928785
$ws.Range("J5").Value = 0.4303743900928785
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.273683333333334
$ws.Range("N5").Value = 18.82105
$ws.Range("O5").Value = 0.5180127793393899
$ws.Range("P5").Value = 0.51801277933939
$ws.Range("Q5").Value = 59.80468522401111
$ws.Range("R5").Value = 538.2421670161
$ws.Range("S5").Value = 0.2229394339685068
$ws.Range("T5").Value = 0.2229394339685068
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 9.532627333333332
$ws.Range("H6").Value = 28.597882
$ws.Range("I6").Value = 0.4303743900928785
$ws.Range("J6").Value = 0.4303743900928785
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.455641666666668
$ws.Range("N6").Value = 16.366925
$ws.Range("O6").Value = 0.4504677639392779
$ws.Range("P6").Value = 0.450467763939278
$ws.Range("Q6").Value = 52.00659887253889
$ws.Range("R6").Value = 468.05938985285
$ws.Range("S6").Value = 0.1938697891618695
$ws.Range("T6").Value = 0.1938697891618695
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.532627333333332
$ws.Range("H7").Value = 28.597882
$ws.Range("I7").Value = 0.4303743900928785
$ws.Range("J7").Value = 0.4303743900928785
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.381734
$ws.Range("N7").Value = 1.145202
$ws.Range("O7").Value = 0.03151945672133213
$ws.Range("P7").Value = 0.03151945672133213
$ws.Range("Q7").Value = 3.638927962462666
$ws.Range("R7").Value = 32.75035166216399
$ws.Range("S7").Value = 0.01356516696250219
$ws.Range("T7").Value = 0.01356516696250219
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.524702666666666
$ws.Range("H8").Value = 10.574108
$ws.Range("I8").Value = 0.1591315497167317
$ws.Range("J8").Value = 0.1591315497167317
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.273683333333334
$ws.Range("N8").Value = 18.82105
$ws.Range("O8").Value = 0.5180127793393899
$ws.Range("P8").Value = 0.51801277933939
$ws.Range("Q8").Value = 22.11286837482222
$ws.Range("R8").Value = 199.0158153734
$ws.Range("S8").Value = 0.0824321763493485
$ws.Range("T8").Value = 0.08243217634934852
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.524702666666666
$ws.Range("H9").Value = 10.574108
$ws.Range("I9").Value = 0.1591315497167317
$ws.Range("J9").Value = 0.1591315497167317
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.455641666666668
$ws.Range("N9").Value = 16.366925
$ws.Range("O9").Value = 0.4504677639392779
$ws.Range("P9").Value = 0.450467763939278
$ws.Range("Q9").Value = 19.22951473087778
$ws.Range("R9").Value = 173.0656325779
$ws.Range("S9").Value = 0.07168363337308817
$ws.Range("T9").Value = 0.07168363337308817
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.524702666666666
$ws.Range("H10").Value = 10.574108
$ws.Range("I10").Value = 0.1591315497167317
$ws.Range("J10").Value = 0.1591315497167317
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.381734
$ws.Range("N10").Value = 1.145202
$ws.Range("O10").Value = 0.03151945672133213
$ws.Range("P10").Value = 0.03151945672133213
$ws.Range("Q10").Value = 1.345498847757333
$ws.Range("R10").Value = 12.109489629816
$ws.Range("S10").Value = 0.005015739994295037
$ws.Range("T10").Value = 0.005015739994295037
